$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.362.99'
$ws.Range('E2').Value = '  +0.06%  '

$ws.Range('D3').Value = '2.282.88'
$ws.Range('E3').Value = '  -0.40%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '504.32'
$ws.Range('E5').Value = '  +1.69%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.63'
$ws.Range('E6').Value = '  +1.86%  '

$ws.Range('E7').Value = '  -0.23%  '

$ws.Range('E8').Value = '  +0.07%  '

$ws.Range('E9').Value = '  +1.59%  '

$ws.Range('E10').Value = '  +0.98%  '

$ws.Range('E11').Value = '  +4.03%  '

$ws.Range('E12').Value = '  +2.09%  '

$ws.Range('D13').Value = '2.689.74'
$ws.Range('E13').Value = '  -0.26%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.94'
$ws.Range('E14').Value = '  +5.86%  '

$ws.Range('D15').Value = '54.324.19'
$ws.Range('E15').Value = '  +0.12%  '

$ws.Range('E16').Value = '  +0.22%  '

$ws.Range('D17').Value = '2.283.35'
$ws.Range('E17').Value = '  +0.23%  '

$ws.Range('E18').Value = '  +3.86%  '

$ws.Range('E19').Value = '  +2.26%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '305.65'
$ws.Range('E20').Value = '  +0.75%  '

$ws.Range('E21').Value = '  +0.40%  '

$ws.Range('E22').Value = '  -0.01%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '61.85'
$ws.Range('E23').Value = '  -3.09%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  +0.02%  '

$ws.Range('E25').Value = '  +1.72%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.34'
$ws.Range('E26').Value = '  +2.64%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '174.37'
$ws.Range('E27').Value = '  +5.05%  '

$ws.Range('E28').Value = '  +1.76%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.02'
$ws.Range('E29').Value = '  +2.46%  '

$ws.Range('E30').Value = '  +1.56%  '

$ws.Range('E31').Value = '  +1.92%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '17.82'
$ws.Range('E33').Value = '  +1.04%  '

$ws.Range('E34').Value = '  +10.23%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.997'
$ws.Range('E35').Value = '  -0.19%  '

$ws.Range('E36').Value = '  +1.21%  '

$ws.Range('E37').Value = '  +3.48%  '

$ws.Range('E38').Value = '  -0.34%  '

$ws.Range('E39').Value = '  +0.90%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.40'
$ws.Range('E40').Value = '  +1.48%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.85'
$ws.Range('E41').Value = '  +0.87%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '124.98'
$ws.Range('E42').Value = '  -0.55%  '

$ws.Range('E43').Value = '  +3.40%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0895'
$ws.Range('E44').Value = '  +0.53%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.548'
$ws.Range('E45').Value = '  +0.50%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '243.56'
$ws.Range('E46').Value = '  +2.45%  '

$ws.Range('E47').Value = '  -0.37%  '

$ws.Range('E48').Value = '  +1.13%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.78'
$ws.Range('E49').Value = '  +0.94%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.48'
$ws.Range('E50').Value = '  +0.80%  '

$ws.Range('E51').Value = '  +0.18%  '
